$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.244.41'
$ws.Range('E2').Value = '  +5.23%  '
$ws.Range('D3').Value = '2.765.87'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '580.43'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.13'
$ws.Range('E6').Value = '  +6.84%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.609'
$ws.Range('E8').Value = '  +1.46%  '
$ws.Range('D9').Value = '2.761.94'
$ws.Range('E9').Value = '  +3.34%  '
$ws.Range('E10').Value = '  +1.85%  '
$ws.Range('E11').Value = '  +4.67%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.390'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').Value = '3.251.98'
$ws.Range('E14').Value = '  +3.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.51'
$ws.Range('E15').Value = '  +2.21%  '
$ws.Range('D16').Value = '64.137.75'
$ws.Range('E16').Value = '  +5.07%  '
$ws.Range('E17').Value = '  +6.06%  '
$ws.Range('D18').Value = '2.761.67'
$ws.Range('E18').Value = '  +3.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.02'
$ws.Range('E19').Value = '  +2.93%  '
$ws.Range('E20').Value = '  +2.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '361.47'
$ws.Range('E21').Value = '  +2.90%  '
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.534'
$ws.Range('E24').Value = '  -0.34%  '
$ws.Range('E26').Value = '  +5.55%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.54'
$ws.Range('E27').Value = '  +4.64%  '
$ws.Range('E28').Value = '  -0.51%  '
$ws.Range('D29').Value = '0.0₃0913'
$ws.Range('E29').Value = '  +11.81%  '
$ws.Range('E30').Value = '  +2.27%  '
$ws.Range('E31').Value = '  +3.10%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.27'
$ws.Range('E32').Value = '  +17.79%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '171.96'
$ws.Range('E33').Value = '  +3.90%  '
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '20.49'
$ws.Range('E35').Value = '  +2.59%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.83'
$ws.Range('E36').Value = '  +7.50%  '
$ws.Range('E37').Value = '  +8.47%  '
$ws.Range('E38').Value = '  +9.24%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.02'
$ws.Range('E39').Value = '  +14.78%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '348.08'
$ws.Range('E40').Value = '  +3.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.26'
$ws.Range('E41').Value = '  +5.63%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.19'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.74'
$ws.Range('E43').Value = '  +10.60%  '
$ws.Range('E44').Value = '  +6.66%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.77'
$ws.Range('E45').Value = '  +6.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0590'
$ws.Range('E46').Value = '  +4.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.647'
$ws.Range('E47').Value = '  +5.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '137.61'
$ws.Range('E48').Value = '  +2.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0255'
$ws.Range('E49').Value = '  +2.39%  '
$ws.Range('E50').Value = '  +0.85%  '
